# Add "Wins", "Losses", "Ties" columns (AC, AD, AE) to the worksheet,
# populating each player's season record (rows 2-50) with 81 wins,
# 81 losses, and 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the formatting of the other header cells (bold, bordered, centered)
# by copying the format from an existing header cell.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Data rows 2-50: Wins = 81, Losses = 81, Ties = 0
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 29).Value = 81   # AC
    $ws.Cells.Item($row, 30).Value = 81   # AD
    $ws.Cells.Item($row, 31).Value = 0    # AE
}
